# Auto-generated: apply updated cryptocurrency price/volume data
# (coinranking.com scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.957.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "'1.884.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'305.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.5142"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("D8").Value = "'0.3736"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.56%  "
$ws.Range("D9").Value = "'0.07185"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'21.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'0.8982"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "'0.07645"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "'1.879.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'93.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'5.229"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'0.000008481"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'14.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'27.008.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'5.042"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'2.131.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'6.378"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'146.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.278"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.66%  "
$ws.Range("D27").Value = "'18.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'1.723"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "'113.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'4.897"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.88%  "
$ws.Range("D31").Value = "'4.767"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").Value = "'0.09173"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'0.05030"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'1.231"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.97%  "
$ws.Range("D35").Value = "'0.7635"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "'2.975"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'2.584"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'0.5583"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'0.01986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'9.043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.94%  "
$ws.Range("D43").Value = "'6.606"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "'118.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "'0.1498"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "'0.4804"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'1.591"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "'37.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "'63.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
